# Update Name of Algo
# Apply updated KNN-imputed values for the terrestrial_mammals AC/10/seed3 result set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = -12.526
$ws.Range("A10").Value = -21.736
$ws.Range("A12").Value = -21.696
$ws.Range("C15").Value = -13.895
$ws.Range("A18").Value = -22.127
$ws.Range("C20").Value = -12.46
$ws.Range("C29").Value = -12.124
$ws.Range("C30").Value = -13.347
$ws.Range("C31").Value = -13.358
$ws.Range("A37").Value = -20.029
$ws.Range("C40").Value = -12.782
$ws.Range("A55").Value = -21.868
$ws.Range("A68").Value = -21.736
$ws.Range("C68").Value = -11.001
$ws.Range("C76").Value = -13.045
$ws.Range("A77").Value = -20.843
$ws.Range("A78").Value = -19.951
$ws.Range("C87").Value = -13.458
$ws.Range("C88").Value = -13.35
$ws.Range("C96").Value = -12.586
$ws.Range("C98").Value = -13.23
$ws.Range("C101").Value = -13.048
$ws.Range("C102").Value = -13.091
